$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.182636380195618
$ws.Range("B1").Value = 2.157172203063965
$ws.Range("C1").Value = 3.012982606887817
$ws.Range("D1").Value = 3.518430471420288
$ws.Range("E1").Value = 1.633718490600586
